$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set formulas for G11 and H11 (new data for Slide 3)
$ws.Range("G11").Formula = "=4000/60"
$ws.Range("H11").Formula = "=G11-G10"

# Update the active selection on the sheet view
$ws.Range("E19").Select()

$wb.Save()
